$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 11:23"

# Row 25 - Filipinas
$ws.Range("B25").Value = 252964
$ws.Range("C25").Value = 4040
$ws.Range("D25").Value = 186606
$ws.Range("E25").Value = 62250
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = 4108

# Row 26 - Indonesia
$ws.Range("B26").Value = 210940
$ws.Range("C26").Value = 3737
$ws.Range("D26").Value = 150217
$ws.Range("E26").Value = 52179
$ws.Range("G26").Value = 88
$ws.Range("H26").Value = 8544

# Row 49 - Polonia
$ws.Range("B49").Value = 73047
$ws.Range("C49").Value = 594
$ws.Range("D49").Value = 58848
$ws.Range("E49").Value = 12030
$ws.Range("G49").Value = 10
$ws.Range("H49").Value = 2169

# Row 55 - Singapur
$ws.Range("B55").Value = 57316
$ws.Range("C55").Value = 87
$ws.Range("E55").Value = 731

# Row 90 - Croacia
$ws.Range("B90").Value = 13107
$ws.Range("C90").Value = 190
$ws.Range("D90").Value = 10466
$ws.Range("E90").Value = 2430
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 211

# Row 103 - Finlandia
$ws.Range("B103").Value = 8512
$ws.Range("C103").Value = 43
$ws.Range("E103").Value = 675

# Row 112 - Eslovaquia
$ws.Range("B112").Value = 5252
$ws.Range("C112").Value = 186
$ws.Range("D112").Value = 3081
$ws.Range("E112").Value = 2134

# Row 115 - Hong Kong
$ws.Range("B115").Value = 4926
$ws.Range("C115").Value = 12
$ws.Range("D115").Value = 4597
$ws.Range("E115").Value = 230

# Row 159 - Letonia
$ws.Range("B159").Value = 1459
$ws.Range("C159").Value = 11
$ws.Range("E159").Value = 176

# Row 176 - Taiwan
$ws.Range("B176").Value = 498
$ws.Range("C176").Value = 2
$ws.Range("E176").Value = 16

$wb.Save()
